$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 106/107, pushing the existing data (old rows 106-159)
# down to rows 108-161.
$ws.Rows("106:107").Insert()

# New row 106: Zafiro rojo entry dated 2021-10-20 (serial 44489)
$ws.Range("A106").Value = 7
$ws.Range("B106").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C106").Value = 'Ñuble'
$ws.Range("D106").Value = 44489
$ws.Range("E106").Value = 16
$ws.Range("F106").Value = 100112002
$ws.Range("G106").Value = 'Pimiento'
$ws.Range("H106").Value = 'Zafiro rojo'
$ws.Range("I106").Value = 'Primera'
$ws.Range("J106").Value = 120
$ws.Range("K106").Value = 43000
$ws.Range("L106").Value = 44000
$ws.Range("M106").Value = 43500
$ws.Range("N106").Value = '$/caja 15 kilos'
$ws.Range("O106").Value = 'Región de Arica y Parinacota'
$ws.Range("P106").Value = 2900
$ws.Range("Q106").Value = 15
$ws.Range("R106").Value = 'Hortaliza'

# New row 107: Zafiro verde entry dated 2021-10-20 (serial 44489)
$ws.Range("A107").Value = 7
$ws.Range("B107").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C107").Value = 'Ñuble'
$ws.Range("D107").Value = 44489
$ws.Range("E107").Value = 16
$ws.Range("F107").Value = 100112002
$ws.Range("G107").Value = 'Pimiento'
$ws.Range("H107").Value = 'Zafiro verde'
$ws.Range("I107").Value = 'Primera'
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 41000
$ws.Range("L107").Value = 42000
$ws.Range("M107").Value = 41500
$ws.Range("N107").Value = '$/caja 15 kilos'
$ws.Range("O107").Value = 'Región de Arica y Parinacota'
$ws.Range("P107").Value = 2767
$ws.Range("Q107").Value = 15
$ws.Range("R107").Value = 'Hortaliza'
